$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" date; update every data row (2-89)
# from serial date 45243 (2023-11-13) to 45244 (2023-11-14).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
